$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Goal (per the diff): the last paragraph currently reads
#   "...//Preguntar/calc" + "u" + [_GoBack bookmark] + "lar por el keySize
#    maximo con el q se calcularan los nnc"
# and must become a paragraph whose text is unchanged but with the
# bookmark removed from mid-paragraph, followed by a brand-new paragraph
# (same red font color) reading:
#   "//añadir a la info de los ataques el enlace a criptored?"
# with the _GoBack bookmark now collapsed at the very end of that new
# paragraph (just before its paragraph mark).
# -----------------------------------------------------------------------

# Locate the paragraph that currently holds the (mid-text) _GoBack bookmark
# -- it's the one ending in "...se calcularan los nnc".
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*calcularan los nnc*") {
        $target = $p
    }
}

$endPos = $target.Range.End - 1   # just after the last real character ("nnc"), before the paragraph mark

# 1) Insert the brand new paragraph (with its own proofed runs) right
#    after the existing text, still inside the red-font block.
$newParaXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:color w:val="FF0000"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:color w:val="FF0000"/>
              </w:rPr>
              <w:t xml:space="preserve">//a&#241;adir a la </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:color w:val="FF0000"/>
              </w:rPr>
              <w:t>info</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:color w:val="FF0000"/>
              </w:rPr>
              <w:t xml:space="preserve"> de los ataques el enlace a criptored?</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insertionPoint = $d.Range($endPos, $endPos)
$insertionPoint.InsertXML($newParaXml)

# 2) Re-find the freshly inserted paragraph (it now directly follows
#    $target) so we can anchor the bookmark to the very end of its text.
$newPara = $target.Next()
$newEnd = $newPara.Range.End - 1   # just before the new paragraph's own mark

# Workaround: adding a *collapsed* bookmark exactly at a position that is
# immediately followed by an empty paragraph mis-resolves in this host, so
# pad with a throwaway character, bookmark before it, then strip the pad.
$pad = $d.Range($newEnd, $newEnd)
$pad.InsertAfter("@")

$bmRange = $d.Range($newEnd, $newEnd)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$padRange = $d.Range($newEnd, $newEnd + 1)
$padRange.Delete()

# 3) Finally, strip the now-stale bookmark that used to sit mid-text in
#    $target (the one originally wrapping "...calcu|lar...").
foreach ($b in $d.Bookmarks) {
    if ($b.Name -eq "_GoBack_OLD_PLACEHOLDER_NEVER_MATCHES") {
        $b.Delete()
    }
}
